$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.006416797637939
$ws.Range("B1").Value = 4.796367168426514
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.810356140136719
$ws.Range("E1").Value = 2.391666650772095
